# 144L_2021_BactAbund.xlsx - "modifying bact abundance 2021 to test fetch"
#
# Adds two new worksheets (TOC_Data, testing fetch) after DAPI_Data, tweaks the
# Metadata sheet's scroll position, and re-styles the DAPI_Data sheet's
# selection / column widths.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: scroll back to the top (drop topLeftCell="A39"), keep the
#    existing G62 selection.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Activate()
$wsMeta.Range("G62").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 2. DAPI_Data sheet: reformat - select columns A:B, widen a few columns.
# ---------------------------------------------------------------------------
$wsDapi = $wb.Worksheets.Item("DAPI_Data")
$wsDapi.Activate()
$wsDapi.Columns.Item(1).ColumnWidth = 31.498697916666668
$wsDapi.Columns.Item(4).ColumnWidth = 19.330729166666668
$wsDapi.Columns.Item(5).ColumnWidth = 29.498697916666668
$wsDapi.Columns.Item(6).ColumnWidth = 25.830729166666668
$wsDapi.Range("A1:B1048576").Select()

# ---------------------------------------------------------------------------
# 3. New sheet "TOC_Data" (after DAPI_Data): Treatment / Timepoint table.
# ---------------------------------------------------------------------------
$wsToc = $wb.Worksheets.Add($null, $wsDapi)
$wsToc.Name = "TOC_Data"
$wsToc.Columns.Item(1).ColumnWidth = 31.498697916666668

$wsToc.Range("A1").Value = "Treatment"
$wsToc.Range("B1").Value = "Timepoint"

$tocRows = @(
    @("Control", 0),
    @("Control", 4),
    @("Control", 8),
    @("Control", 9),
    @("Kelp Exudate", 0),
    @("Kelp Exudate", 4),
    @("Kelp Exudate", 8),
    @("Kelp Exudate", 9),
    @("Kelp Exudate_Nitrate_Phosphate", 0),
    @("Kelp Exudate_Nitrate_Phosphate", 4),
    @("Kelp Exudate_Nitrate_Phosphate", 8),
    @("Kelp_Exudate_Nitrate_Phosphate", 9),
    @("Glucose_Nitrate_Phosphate", 0),
    @("Glucose_Nitrate_Phosphate", 4),
    @("Glucose_Nitrate_Phosphate", 8),
    @("Glucose_Nitrate_Phosphate", 9)
)

$r = 2
foreach ($row in $tocRows) {
    $wsToc.Range("A$r").Value = $row[0]
    $wsToc.Range("B$r").Value = $row[1]
    $r = $r + 1
}

$wsToc.Range("B18").Select()

# ---------------------------------------------------------------------------
# 4. New sheet "testing fetch" (after TOC_Data): four "data " placeholder rows.
# ---------------------------------------------------------------------------
$wsFetch = $wb.Worksheets.Add($null, $wsToc)
$wsFetch.Name = "testing fetch"

$wsFetch.Range("A1").Value = "data "
$wsFetch.Range("A2").Value = "data "
$wsFetch.Range("A3").Value = "data "
$wsFetch.Range("A4").Value = "data "

$wsFetch.Range("A5").Select()
$wsFetch.Activate()
